$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates (order matches shared-string append order in target)
$ws.Range("D2").Value = "Snowmagdon"
$ws.Range("E2").Value = "Let is snow!"
$ws.Range("H2").Value = "null"
$ws.Range("C2").Value = "weather"
